$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2 = @{ B=1.505614041169197;  C=1.65323645889881;   D=0.1529057820181812; E=0.4998867070740569; G=3.811642989160245 }
    3 = @{ B=0.1554434735375247; C=1.65323645889881;   D=0.7127328510149897; E=0.4998867070740569; G=3.021299490525381 }
    4 = @{ B=3.182878228561681;  C=1.65323645889881;   D=0.7127328510149897; E=0.4998867070740569; G=6.048734245549538 }
    5 = @{ B=1.505614041169197;  C=1.65323645889881;   D=3.082599426703578;  E=0.4998867070740569; G=6.741336633845642 }
    6 = @{ B=1.505614041169197;  C=0.05231270169004087; D=0.7127328510149897; E=6.48142807727062;   G=8.752087671144849 }
    7 = @{ B=0.3464964993005633; C=1.65323645889881;   D=0.7127328510149897; E=6.48142807727062;   G=9.193893886484982 }
    8 = @{ B=0.3464964993005633; C=0.05231270169004087; D=2938.103010863317;  E=246.9852506941017;  G=3185.487070758409 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
